$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: merge A1:I1 and center the title ---
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1:I1").Merge()

# --- "Récuperer les données" block ---
$ws.Range("A3").Value = "Récuperer les données"
$ws.Range("B4").Value = "télécharger les données 1 fois par jours"
$ws.Range("B5").Value = "dézipper les données"
$ws.Range("B6").Value = "lire les fichiers des main leagues"
$ws.Range("B7").Value = "afficher 1 ligne par match"
$ws.Range("B8").Value = "récuperer l'ensemble des statitstiques possible d'un match"
$ws.Range("B9").Value = "L'ensemble des matchs des divisions majeures"

# --- "Ttraiter les données" block ---
$ws.Range("A11").Value = "Ttraiter les données"
$ws.Range("B12").Value = "etablir le classement de la saison "
$ws.Range("B13").Value = "nombre de points"
$ws.Range("B14").Value = "nombre de victoire"
$ws.Range("B15").Value = "nombre de défaite"
$ws.Range("B16").Value = "nombre de matche nul"
$ws.Range("B17").Value = "nombre de matches joués"
$ws.Range("B18").Value = "nombre de but mis par équipe"
$ws.Range("B19").Value = "nombre de but encaissé par équipe"
$ws.Range("B20").Value = "golaverage"
$ws.Range("B21").Value = "nombre de but moyen par match"

# --- "Afficher les données" block ---
$ws.Range("A23").Value = "Afficher les données"
$ws.Range("B24").Value = "Afficher les statistiques par main leagues"

# --- Column B width ---
$ws.Columns("B").ColumnWidth = 49.1640625

# --- View settings ---
$ws.Application.ActiveWindow.Zoom = 125
$ws.Range("B26").Select() | Out-Null
